$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 27.83862766666667
$ws.Range("H2").Value = 83.515883
$ws.Range("I2").Value = 0.04355088691831899
$ws.Range("J2").Value = 0.04355088691831899
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 94.799851
$ws.Range("N2").Value = 284.399553
$ws.Range("O2").Value = 0.1681963571786457
$ws.Range("P2").Value = 0.1681963571786457
$ws.Range("Q2").Value = 2639.097754844478
$ws.Range("R2").Value = 23751.8797936003
$ws.Range("S2").Value = 0.007325100531560388
$ws.Range("T2").Value = 0.007325100531560391

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 27.83862766666667
$ws.Range("H3").Value = 83.515883
$ws.Range("I3").Value = 0.04355088691831899
$ws.Range("J3").Value = 0.04355088691831899
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 145.6413626666667
$ws.Range("N3").Value = 436.924088
$ws.Range("O3").Value = 0.2584006873076977
$ws.Range("P3").Value = 0.2584006873076978
$ws.Range("Q3").Value = 4054.4556681433
$ws.Range("R3").Value = 36490.10101328971
$ws.Range("S3").Value = 0.01125357911255345
$ws.Range("T3").Value = 0.01125357911255345

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 27.83862766666667
$ws.Range("H4").Value = 83.515883
$ws.Range("I4").Value = 0.04355088691831899
$ws.Range("J4").Value = 0.04355088691831899
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 72.77597066666667
$ws.Range("N4").Value = 218.327912
$ws.Range("O4").Value = 0.1291210168281099
$ws.Range("P4").Value = 0.1291210168281099
$ws.Range("Q4").Value = 2025.983150469589
$ws.Range("R4").Value = 18233.84835422629
$ws.Range("S4").Value = 0.005623334802659375
$ws.Range("T4").Value = 0.005623334802659376

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 27.83862766666667
$ws.Range("H5").Value = 83.515883
$ws.Range("I5").Value = 0.04355088691831899
$ws.Range("J5").Value = 0.04355088691831899
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 76.38610333333334
$ws.Range("N5").Value = 229.15831
$ws.Range("O5").Value = 0.1355262079445491
$ws.Range("P5").Value = 0.1355262079445491
$ws.Range("Q5").Value = 2126.484289604192
$ws.Range("R5").Value = 19138.35860643773
$ws.Range("S5").Value = 0.005902286556661642
$ws.Range("T5").Value = 0.005902286556661644

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 27.83862766666667
$ws.Range("H6").Value = 83.515883
$ws.Range("I6").Value = 0.04355088691831899
$ws.Range("J6").Value = 0.04355088691831899
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 58.18688599999999
$ws.Range("N6").Value = 174.560658
$ws.Range("O6").Value = 0.10323668399826
$ws.Range("P6").Value = 0.1032366839982601
$ws.Range("Q6").Value = 1619.843054436779
$ws.Range("R6").Value = 14578.58748993101
$ws.Range("S6").Value = 0.004496049150630455
$ws.Range("T6").Value = 0.004496049150630455

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 27.83862766666667
$ws.Range("H7").Value = 83.515883
$ws.Range("I7").Value = 0.04355088691831899
$ws.Range("J7").Value = 0.04355088691831899
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 115.8358916666667
$ws.Range("N7").Value = 347.5076749999999
$ws.Range("O7").Value = 0.2055190467427377
$ws.Range("P7").Value = 0.2055190467427377
$ws.Range("Q7").Value = 3224.712258544669
$ws.Range("R7").Value = 29022.41032690202
$ws.Range("S7").Value = 0.008950536764253685
$ws.Range("T7").Value = 0.008950536764253687

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 81.69726833333333
$ws.Range("H8").Value = 245.091805
$ws.Range("I8").Value = 0.1278076109685829
$ws.Range("J8").Value = 0.1278076109685829
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 94.799851
$ws.Range("N8").Value = 284.399553
$ws.Range("O8").Value = 0.1681963571786457
$ws.Range("P8").Value = 0.1681963571786457
$ws.Range("Q8").Value = 7744.888865107018
$ws.Range("R8").Value = 69703.99978596317
$ws.Range("S8").Value = 0.02149677458462117
$ws.Range("T8").Value = 0.02149677458462117

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 81.69726833333333
$ws.Range("H9").Value = 245.091805
$ws.Range("I9").Value = 0.1278076109685829
$ws.Range("J9").Value = 0.1278076109685829
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 145.6413626666667
$ws.Range("N9").Value = 436.924088
$ws.Range("O9").Value = 0.2584006873076977
$ws.Range("P9").Value = 0.2584006873076978
$ws.Range("Q9").Value = 11898.50148621098
$ws.Range("R9").Value = 107086.5133758988
$ws.Range("S9").Value = 0.03302557451743667
$ws.Range("T9").Value = 0.03302557451743668

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 81.69726833333333
$ws.Range("H10").Value = 245.091805
$ws.Range("I10").Value = 0.1278076109685829
$ws.Range("J10").Value = 0.1278076109685829
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 72.77597066666667
$ws.Range("N10").Value = 218.327912
$ws.Range("O10").Value = 0.1291210168281099
$ws.Range("P10").Value = 0.1291210168281099
$ws.Range("Q10").Value = 5945.598003773462
$ws.Range("R10").Value = 53510.38203396116
$ws.Range("S10").Value = 0.01650264868663491
$ws.Range("T10").Value = 0.01650264868663491

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 81.69726833333333
$ws.Range("H11").Value = 245.091805
$ws.Range("I11").Value = 0.1278076109685829
$ws.Range("J11").Value = 0.1278076109685829
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 76.38610333333334
$ws.Range("N11").Value = 229.15831
$ws.Range("O11").Value = 0.1355262079445491
$ws.Range("P11").Value = 0.1355262079445491
$ws.Range("Q11").Value = 6240.535980961061
$ws.Range("R11").Value = 56164.82382864955
$ws.Range("S11").Value = 0.0173212808610242
$ws.Range("T11").Value = 0.0173212808610242

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 81.69726833333333
$ws.Range("H12").Value = 245.091805
$ws.Range("I12").Value = 0.1278076109685829
$ws.Range("J12").Value = 0.1278076109685829
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 58.18688599999999
$ws.Range("N12").Value = 174.560658
$ws.Range("O12").Value = 0.10323668399826
$ws.Range("P12").Value = 0.1032366839982601
$ws.Range("Q12").Value = 4753.709639023075
$ws.Range("R12").Value = 42783.38675120768
$ws.Range("S12").Value = 0.01319443394613615
$ws.Range("T12").Value = 0.01319443394613615

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 81.69726833333333
$ws.Range("H13").Value = 245.091805
$ws.Range("I13").Value = 0.1278076109685829
$ws.Range("J13").Value = 0.1278076109685829
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 115.8358916666667
$ws.Range("N13").Value = 347.5076749999999
$ws.Range("O13").Value = 0.2055190467427377
$ws.Range("P13").Value = 0.2055190467427377
$ws.Range("Q13").Value = 9463.475924122595
$ws.Range("R13").Value = 85171.28331710336
$ws.Range("S13").Value = 0.02626689837272983
$ws.Range("T13").Value = 0.02626689837272983

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 269.5867006666667
$ws.Range("H14").Value = 808.7601020000001
$ws.Range("I14").Value = 0.4217427689323495
$ws.Range("J14").Value = 0.4217427689323495
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 94.799851
$ws.Range("N14").Value = 284.399553
$ws.Range("O14").Value = 0.1681963571786457
$ws.Range("P14").Value = 0.1681963571786457
$ws.Range("Q14").Value = 25556.7790547816
$ws.Range("R14").Value = 230011.0114930345
$ws.Range("S14").Value = 0.07093559740085648
$ws.Range("T14").Value = 0.0709355974008565

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 269.5867006666667
$ws.Range("H15").Value = 808.7601020000001
$ws.Range("I15").Value = 0.4217427689323495
$ws.Range("J15").Value = 0.4217427689323495
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 145.6413626666667
$ws.Range("N15").Value = 436.924088
$ws.Range("O15").Value = 0.2584006873076977
$ws.Range("P15").Value = 0.2584006873076978
$ws.Range("Q15").Value = 39262.9744419041
$ws.Range("R15").Value = 353366.769977137
$ws.Range("S15").Value = 0.1089786213591706
$ws.Range("T15").Value = 0.1089786213591707

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 269.5867006666667
$ws.Range("H16").Value = 808.7601020000001
$ws.Range("I16").Value = 0.4217427689323495
$ws.Range("J16").Value = 0.4217427689323495
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 72.77597066666667
$ws.Range("N16").Value = 218.327912
$ws.Range("O16").Value = 0.1291210168281099
$ws.Range("P16").Value = 0.1291210168281099
$ws.Range("Q16").Value = 19619.43381984078
$ws.Range("R16").Value = 176574.904378567
$ws.Range("S16").Value = 0.05445585516444754
$ws.Range("T16").Value = 0.05445585516444755

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 269.5867006666667
$ws.Range("H17").Value = 808.7601020000001
$ws.Range("I17").Value = 0.4217427689323495
$ws.Range("J17").Value = 0.4217427689323495
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 76.38610333333334
$ws.Range("N17").Value = 229.15831
$ws.Range("O17").Value = 0.1355262079445491
$ws.Range("P17").Value = 0.1355262079445491
$ws.Range("Q17").Value = 20592.6775744164
$ws.Range("R17").Value = 185334.0981697477
$ws.Range("S17").Value = 0.05715719820143551
$ws.Range("T17").Value = 0.05715719820143553

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 269.5867006666667
$ws.Range("H18").Value = 808.7601020000001
$ws.Range("I18").Value = 0.4217427689323495
$ws.Range("J18").Value = 0.4217427689323495
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 58.18688599999999
$ws.Range("N18").Value = 174.560658
$ws.Range("O18").Value = 0.10323668399826
$ws.Range("P18").Value = 0.1032366839982601
$ws.Range("Q18").Value = 15686.41061880746
$ws.Range("R18").Value = 141177.6955692671
$ws.Range("S18").Value = 0.04353932496482017
$ws.Range("T18").Value = 0.04353932496482017

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 269.5867006666667
$ws.Range("H19").Value = 808.7601020000001
$ws.Range("I19").Value = 0.4217427689323495
$ws.Range("J19").Value = 0.4217427689323495
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 115.8358916666667
$ws.Range("N19").Value = 347.5076749999999
$ws.Range("O19").Value = 0.2055190467427377
$ws.Range("P19").Value = 0.2055190467427377
$ws.Range("Q19").Value = 31227.81585319809
$ws.Range("R19").Value = 281050.3426787828
$ws.Range("S19").Value = 0.08667617184161916
$ws.Range("T19").Value = 0.08667617184161917

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 218.6560773333333
$ws.Range("H20").Value = 655.9682320000001
$ws.Range("I20").Value = 0.3420666496915519
$ws.Range("J20").Value = 0.3420666496915519
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 94.799851
$ws.Range("N20").Value = 284.399553
$ws.Range("O20").Value = 0.1681963571786457
$ws.Range("P20").Value = 0.1681963571786457
$ws.Range("Q20").Value = 20728.56355144448
$ws.Range("R20").Value = 186557.0719630003
$ws.Range("S20").Value = 0.05753436439042293
$ws.Range("T20").Value = 0.05753436439042294

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 218.6560773333333
$ws.Range("H21").Value = 655.9682320000001
$ws.Range("I21").Value = 0.3420666496915519
$ws.Range("J21").Value = 0.3420666496915519
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 145.6413626666667
$ws.Range("N21").Value = 436.924088
$ws.Range("O21").Value = 0.2584006873076977
$ws.Range("P21").Value = 0.2584006873076978
$ws.Range("Q21").Value = 31845.36905817471
$ws.Range("R21").Value = 286608.3215235724
$ws.Range("S21").Value = 0.08839025738533848
$ws.Range("T21").Value = 0.0883902573853385

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 218.6560773333333
$ws.Range("H22").Value = 655.9682320000001
$ws.Range("I22").Value = 0.3420666496915519
$ws.Range("J22").Value = 0.3420666496915519
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 72.77597066666667
$ws.Range("N22").Value = 218.327912
$ws.Range("O22").Value = 0.1291210168281099
$ws.Range("P22").Value = 0.1291210168281099
$ws.Range("Q22").Value = 15912.90827009907
$ws.Range("R22").Value = 143216.1744308916
$ws.Range("S22").Value = 0.04416799363115804
$ws.Range("T22").Value = 0.04416799363115804

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 218.6560773333333
$ws.Range("H23").Value = 655.9682320000001
$ws.Range("I23").Value = 0.3420666496915519
$ws.Range("J23").Value = 0.3420666496915519
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 76.38610333333334
$ws.Range("N23").Value = 229.15831
$ws.Range("O23").Value = 0.1355262079445491
$ws.Range("P23").Value = 0.1355262079445491
$ws.Range("Q23").Value = 16702.28571764533
$ws.Range("R23").Value = 150320.5714588079
$ws.Range("S23").Value = 0.04635899589699249
$ws.Range("T23").Value = 0.0463589958969925

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 218.6560773333333
$ws.Range("H24").Value = 655.9682320000001
$ws.Range("I24").Value = 0.3420666496915519
$ws.Range("J24").Value = 0.3420666496915519
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 58.18688599999999
$ws.Range("N24").Value = 174.560658
$ws.Range("O24").Value = 0.10323668399826
$ws.Range("P24").Value = 0.1032366839982601
$ws.Range("Q24").Value = 12722.91624500185
$ws.Range("R24").Value = 114506.2462050167
$ws.Range("S24").Value = 0.03531382662055026
$ws.Range("T24").Value = 0.03531382662055026

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 218.6560773333333
$ws.Range("H25").Value = 655.9682320000001
$ws.Range("I25").Value = 0.3420666496915519
$ws.Range("J25").Value = 0.3420666496915519
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 115.8358916666667
$ws.Range("N25").Value = 347.5076749999999
$ws.Range("O25").Value = 0.2055190467427377
$ws.Range("P25").Value = 0.2055190467427377
$ws.Range("Q25").Value = 25328.22168624229
$ws.Range("R25").Value = 227953.9951761806
$ws.Range("S25").Value = 0.07030121176708976
$ws.Range("T25").Value = 0.07030121176708976

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 30.15805666666667
$ws.Range("H26").Value = 90.47417000000002
$ws.Range("I26").Value = 0.04717941312670751
$ws.Range("J26").Value = 0.04717941312670752
$ws.Range("K26").Value = 3
$ws.Range("M26").Value = 94.799851
$ws.Range("N26").Value = 284.399553
$ws.Range("O26").Value = 0.1681963571786457
$ws.Range("P26").Value = 0.1681963571786457
$ws.Range("Q26").Value = 2858.979278449557
$ws.Range("R26").Value = 25730.81350604602
$ws.Range("S26").Value = 0.007935405421738581
$ws.Range("T26").Value = 0.007935405421738583

$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 30.15805666666667
$ws.Range("H27").Value = 90.47417000000002
$ws.Range("I27").Value = 0.04717941312670751
$ws.Range("J27").Value = 0.04717941312670752
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 145.6413626666667
$ws.Range("N27").Value = 436.924088
$ws.Range("O27").Value = 0.2584006873076977
$ws.Range("P27").Value = 0.2584006873076978
$ws.Range("Q27").Value = 4392.260468311884
$ws.Range("R27").Value = 39530.34421480697
$ws.Range("S27").Value = 0.01219119277871504
$ws.Range("T27").Value = 0.01219119277871504

$ws.Range("E28").Value = 3
$ws.Range("G28").Value = 30.15805666666667
$ws.Range("H28").Value = 90.47417000000002
$ws.Range("I28").Value = 0.04717941312670751
$ws.Range("J28").Value = 0.04717941312670752
$ws.Range("K28").Value = 3
$ws.Range("M28").Value = 72.77597066666667
$ws.Range("N28").Value = 218.327912
$ws.Range("O28").Value = 0.1291210168281099
$ws.Range("P28").Value = 0.1291210168281099
$ws.Range("Q28").Value = 2194.781847337005
$ws.Range("R28").Value = 19753.03662603304
$ws.Range("S28").Value = 0.006091853796273947
$ws.Range("T28").Value = 0.006091853796273948

$ws.Range("E29").Value = 3
$ws.Range("G29").Value = 30.15805666666667
$ws.Range("H29").Value = 90.47417000000002
$ws.Range("I29").Value = 0.04717941312670751
$ws.Range("J29").Value = 0.04717941312670752
$ws.Range("K29").Value = 3
$ws.Range("M29").Value = 76.38610333333334
$ws.Range("N29").Value = 229.15831
$ws.Range("O29").Value = 0.1355262079445491
$ws.Range("P29").Value = 0.1355262079445491
$ws.Range("Q29").Value = 2303.656432872523
$ws.Range("R29").Value = 20732.90789585271
$ws.Range("S29").Value = 0.00639404695411195
$ws.Range("T29").Value = 0.006394046954111953

$ws.Range("E30").Value = 3
$ws.Range("G30").Value = 30.15805666666667
$ws.Range("H30").Value = 90.47417000000002
$ws.Range("I30").Value = 0.04717941312670751
$ws.Range("J30").Value = 0.04717941312670752
$ws.Range("K30").Value = 3
$ws.Range("M30").Value = 58.18688599999999
$ws.Range("N30").Value = 174.560658
$ws.Range("O30").Value = 0.10323668399826
$ws.Range("P30").Value = 0.1032366839982601
$ws.Range("Q30").Value = 1754.803405244873
$ws.Range("R30").Value = 15793.23064720386
$ws.Range("S30").Value = 0.004870646164185265
$ws.Range("T30").Value = 0.004870646164185266

$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 30.15805666666667
$ws.Range("H31").Value = 90.47417000000002
$ws.Range("I31").Value = 0.04717941312670751
$ws.Range("J31").Value = 0.04717941312670752
$ws.Range("K31").Value = 3
$ws.Range("M31").Value = 115.8358916666667
$ws.Range("N31").Value = 347.5076749999999
$ws.Range("O31").Value = 0.2055190467427377
$ws.Range("P31").Value = 0.2055190467427377
$ws.Range("Q31").Value = 3493.385384917195
$ws.Range("R31").Value = 31440.46846425475
$ws.Range("S31").Value = 0.009696268011682734
$ws.Range("T31").Value = 0.009696268011682736

$ws.Range("E32").Value = 3
$ws.Range("G32").Value = 11.283952
$ws.Range("H32").Value = 33.851856
$ws.Range("I32").Value = 0.01765267036248923
$ws.Range("J32").Value = 0.01765267036248923
$ws.Range("K32").Value = 3
$ws.Range("M32").Value = 94.799851
$ws.Range("N32").Value = 284.399553
$ws.Range("O32").Value = 0.1681963571786457
$ws.Range("P32").Value = 0.1681963571786457
$ws.Range("Q32").Value = 1069.716968291152
$ws.Range("R32").Value = 9627.452714620369
$ws.Range("S32").Value = 0.002969114849446131
$ws.Range("T32").Value = 0.002969114849446132

$ws.Range("E33").Value = 3
$ws.Range("G33").Value = 11.283952
$ws.Range("H33").Value = 33.851856
$ws.Range("I33").Value = 0.01765267036248923
$ws.Range("J33").Value = 0.01765267036248923
$ws.Range("K33").Value = 3
$ws.Range("M33").Value = 145.6413626666667
$ws.Range("N33").Value = 436.924088
$ws.Range("O33").Value = 0.2584006873076977
$ws.Range("P33").Value = 0.2584006873076978
$ws.Range("Q33").Value = 1643.410145545259
$ws.Range("R33").Value = 14790.69130990733
$ws.Range("S33").Value = 0.004561462154483442
$ws.Range("T33").Value = 0.004561462154483444

$ws.Range("E34").Value = 3
$ws.Range("G34").Value = 11.283952
$ws.Range("H34").Value = 33.851856
$ws.Range("I34").Value = 0.01765267036248923
$ws.Range("J34").Value = 0.01765267036248923
$ws.Range("K34").Value = 3
$ws.Range("M34").Value = 72.77597066666667
$ws.Range("N34").Value = 218.327912
$ws.Range("O34").Value = 0.1291210168281099
$ws.Range("P34").Value = 0.1291210168281099
$ws.Range("Q34").Value = 821.2005597560747
$ws.Range("R34").Value = 7390.805037804672
$ws.Range("S34").Value = 0.002279330746936048
$ws.Range("T34").Value = 0.002279330746936048

$ws.Range("E35").Value = 3
$ws.Range("G35").Value = 11.283952
$ws.Range("H35").Value = 33.851856
$ws.Range("I35").Value = 0.01765267036248923
$ws.Range("J35").Value = 0.01765267036248923
$ws.Range("K35").Value = 3
$ws.Range("M35").Value = 76.38610333333334
$ws.Range("N35").Value = 229.15831
$ws.Range("O35").Value = 0.1355262079445491
$ws.Range("P35").Value = 0.1355262079445491
$ws.Range("Q35").Value = 861.9371234803733
$ws.Range("R35").Value = 7757.43411132336
$ws.Range("S35").Value = 0.002392399474323294
$ws.Range("T35").Value = 0.002392399474323295

$ws.Range("E36").Value = 3
$ws.Range("G36").Value = 11.283952
$ws.Range("H36").Value = 33.851856
$ws.Range("I36").Value = 0.01765267036248923
$ws.Range("J36").Value = 0.01765267036248923
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 58.18688599999999
$ws.Range("N36").Value = 174.560658
$ws.Range("O36").Value = 0.10323668399826
$ws.Range("P36").Value = 0.1032366839982601
$ws.Range("Q36").Value = 656.5780286534718
$ws.Range("R36").Value = 5909.202257881248
$ws.Range("S36").Value = 0.001822403151937751
$ws.Range("T36").Value = 0.001822403151937752

$ws.Range("E37").Value = 3
$ws.Range("G37").Value = 11.283952
$ws.Range("H37").Value = 33.851856
$ws.Range("I37").Value = 0.01765267036248923
$ws.Range("J37").Value = 0.01765267036248923
$ws.Range("K37").Value = 3
$ws.Range("M37").Value = 115.8358916666667
$ws.Range("N37").Value = 347.5076749999999
$ws.Range("O37").Value = 0.2055190467427377
$ws.Range("P37").Value = 0.2055190467427377
$ws.Range("Q37").Value = 1307.086641443866
$ws.Range("R37").Value = 11763.7797729948
$ws.Range("S37").Value = 0.003627959985362565
$ws.Range("T37").Value = 0.003627959985362565
